$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 3 new rows before row 32 (pushes old row 32 and everything
#        below it down by 3; formulas / merged cells shift automatically) ---
$ws.Rows("32:34").Insert()

# The 3 freshly inserted rows (32:34) don't inherit the right "interior of
# block" borders/number formats. Copy them from neighbouring cells that
# already carry the exact look we need: B from B31 (date column, blank,
# bordered), E from E31 (time column, bordered on the right), F from F83
# (task column, wraps text, bordered on the right only - this is the style
# used for "interior" rows near the bottom of a work block).
$ws.Range("B31").Copy()
$ws.Range("B32:B34").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("E31").Copy()
$ws.Range("E32:E34").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F83").Copy()
$ws.Range("F32:F34").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 2. Fill in the values for the expanded work block (rows 30-34) and
#        its total row (35) ---

# Row 30 (date 27.04.2017 already present) - fill time range + task text
$ws.Range("C30").Value = 0.33333333333333331
$ws.Range("D30").Value = "-"
$ws.Range("E30").Value = 0.44791666666666669
$ws.Range("F30").Value = "Vue editShoot crée et intégrée / rédaction rapport"

# Row 31
$ws.Range("C31").Value = 0.44791666666666669
$ws.Range("D31").Value = "-"
$ws.Range("E31").Value = 0.46875
$ws.Range("F31").Value = "Visite de l'expert 2 M.Malherbe"

# Row 32 (new)
$ws.Range("C32").Value = 0.46875
$ws.Range("D32").Value = "-"
$ws.Range("E32").Value = 0.5
$ws.Range("F32").Value = "Mise au point rapport"

# Row 33 (new)
$ws.Range("C33").Value = 0.53125
$ws.Range("D33").Value = "-"
$ws.Range("E33").Value = 0.64583333333333337
$ws.Range("F33").Value = "Vue editShoot débutée"

# Row 34 (new)
$ws.Range("C34").Value = 0.64583333333333337
$ws.Range("D34").Value = "-"
$ws.Range("E34").Value = 0.71458333333333324
$ws.Range("F34").Value = "Discussion avec chef de projet thèmes / selecteur d'Arrow pour editShoot"

# Row 35 - total row for this block, now covers C30:E34
$ws.Range("E35").Formula = "=E30-C30+E31-C31+E32-C32+E33-C33+E34-C34"

# --- 3. Update the view: selection moved to E35 ---
$ws.Range("E35").Select()

# --- 4. Print area grows by 3 rows (was A1:I82) ---
$ws.PageSetup.PrintArea = "A1:I85"

$wb.Save()
